$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fix a typo in the existing data (row 4, Author column): "SaLma Khaled" -> "Salma Khaled"
$ws.Range("D4").Value = "Salma Khaled"

# --- Append three new review rows (6, 7, 8), matching the formatting of row 5 ---
$ws.Range("A5:G5").Copy()
$ws.Range("A6:G8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 6: SRS review
$ws.Range("A6").Value = "OMS-REV-SRS-05"
$ws.Range("B6").Value = "Requirement/OnlineMobileStore_SRS"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = "M. Ramzy, M. Kassas"
$ws.Range("E6").Value = "Salma Khaled"
$ws.Range("F6").Value = 45539
$comments6 = "1-The table of features is not complete`n2-Should be called as features not modules`n3-User ID in req number ""OMS_SRS_Adm_UD_01"" cause a conflict and should be deleted`n4-The req number ""OMS_SRS_Core-B_04"" is missing according to the naming convention`n5-Req ""OMS_SRS_Reg_01"", you should refer to an existing feature at the last sentence ""redirect them to the login page"" `n6-Req number ""OMS_SRS_Log_01"" is not clear, remove the first sentence`n7-Conflict between client dashpord/ client homepage, The name must be circulated on all requirements"
$ws.Range("G6").Value = $comments6
$ws.Rows.Item(6).RowHeight = 222

# Row 7: PMP re-review
$ws.Range("A7").Value = "OMS-REV-PMP-06"
$ws.Range("B7").Value = "OnlineMobileStore_PMP"
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = "Salma Mohamed"
$ws.Range("E7").Value = "Salma Khaled"
# F7 uses the plain center/middle style (not the date-numFmt one the rest of column F uses)
$ws.Range("A7").Copy()
$ws.Range("F7").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("F7").Value = "_"
$ws.Range("G7").Value = "No comments"
$ws.Rows.Item(7).RowHeight = 30

# Row 8: CIL re-review
$ws.Range("A8").Value = "OMS-REV-CIL-07"
$ws.Range("B8").Value = "PM/OnlineMobileStore_CIL"
$ws.Range("C8").Value = 2
$ws.Range("D8").Value = "Sama Wagdy"
$ws.Range("E8").Value = "Salma Khaled"
$ws.Range("F8").Value = "_"
$ws.Range("G8").Value = "No comments"
$ws.Rows.Item(8).RowHeight = 30

# --- Column widths widened slightly to fit the new (longer) content ---
$ws.Columns.Item(2).ColumnWidth = 34.66796875
$ws.Columns.Item(5).ColumnWidth = 18.65234375
$ws.Columns.Item(6).ColumnWidth = 14.453125
$ws.Columns.Item(7).ColumnWidth = 44.921875

# --- View tweaks to match the saved selection/scroll state ---
$ws.Range("E4").Select()
$ws.Application.ActiveWindow.ScrollRow = 4
